$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "average" (column C) value, same for all data rows 2-9
$ws.Range("C2:C9").Value = 85.67030377800782

# New "p" (column D) values, row-by-row
$ws.Range("D2").Value = 0.99862736463546753
$ws.Range("D3").Value = 0.99875873327255249
$ws.Range("D4").Value = 0.99880015850067139
$ws.Range("D5").Value = 0.99887210130691528
$ws.Range("D6").Value = 0.99896806478500366
$ws.Range("D7").Value = 0.99936532974243164
$ws.Range("D8").Value = 0.99968570470809937
$ws.Range("D9").Value = 0.99984282255172729
